$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gens")
$ws.Range("C2").Value = 1519800.3
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = -1697951.1
$ws.Range("C12").Value = 178985.91
$ws.Range("C22").Value = 0
$ws.Range("C24").Value = 26641.929
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = -27011.879
$ws.Range("C28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("C31").Value = 0

$ws = $wb.Worksheets.Item("lines")
$ws.Range("C2").Value = 1520984.1
$ws.Range("D2").Value = -25
$ws.Range("C3").Value = -7.0890839
$ws.Range("D3").Value = -42116.202
$ws.Range("C4").Value = -1112.763
$ws.Range("D4").Value = 2092.4494
$ws.Range("C5").Value = -92146.077
$ws.Range("D5").Value = -2377.1266
$ws.Range("C6").Value = -84917.936
$ws.Range("D6").Value = 3250.9091
$ws.Range("C7").Value = 15.410916
$ws.Range("D7").Value = 37767.452
$ws.Range("C8").Value = -360
$ws.Range("D8").Value = 291372.53
$ws.Range("E8").Value = 431405.8
$ws.Range("C9").Value = -92220.077
$ws.Range("D9").Value = -1946.6234
$ws.Range("C10").Value = -1183.763
$ws.Range("D10").Value = 2166.3006
$ws.Range("C11").Value = -85053.936
$ws.Range("D11").Value = 1032.8409
$ws.Range("C12").Value = 179060.91
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = 90048.268
$ws.Range("D13").Value = -4303.75
$ws.Range("C14").Value = 88841.637
$ws.Range("D14").Value = 4303.75
$ws.Range("C15").Value = -1119.5241
$ws.Range("D15").Value = 16469.337
$ws.Range("C16").Value = -1211.8738
$ws.Range("D16").Value = 6426.7677
$ws.Range("C17").Value = 1250.6438
$ws.Range("D17").Value = 7861.8374
$ws.Range("C18").Value = 1158.2941
$ws.Range("D18").Value = -2180.7323
$ws.Range("C19").Value = -222.73731
$ws.Range("D19").Value = -6565.9075
$ws.Range("C20").Value = 353.857
$ws.Range("D20").Value = 17910.756
$ws.Range("C21").Value = -61.125359
$ws.Range("D21").Value = 3476.6622
$ws.Range("C22").Value = 7.5456767
$ws.Range("D22").Value = -2122.5949
$ws.Range("C23").Value = 42.137332
$ws.Range("D23").Value = -5599.2571
$ws.Range("C24").Value = 159.857
$ws.Range("D24").Value = 25160.348
$ws.Range("C25").Value = -1114.857
$ws.Range("D25").Value = 56967.685
$ws.Range("C26").Value = 351.4752
$ws.Range("D26").Value = 1186.4407
$ws.Range("C27").Value = 351.4752
$ws.Range("D27").Value = 1186.4407
$ws.Range("C28").Value = 360
$ws.Range("D28").Value = -180373.47
$ws.Range("D29").Value = 95972.374
$ws.Range("E29").Value = 150622.63
$ws.Range("C30").Value = -450
$ws.Range("D30").Value = 1983.5591
$ws.Range("E30").Value = 59893.321
$ws.Range("C31").Value = -10894.119
$ws.Range("D31").Value = -586.44068
$ws.Range("C32").Value = 10894.119
$ws.Range("D32").Value = -686.44068
$ws.Range("C33").Value = 7707.4048
$ws.Range("D33").Value = -544.55206
$ws.Range("C34").Value = 7707.4048
$ws.Range("D34").Value = -544.55206
$ws.Range("C35").Value = -315.5
$ws.Range("D35").Value = 1724.834
$ws.Range("C36").Value = -315.5
$ws.Range("D36").Value = 1724.834
$ws.Range("C37").Value = -379.5
$ws.Range("D37").Value = 948.65868
$ws.Range("C38").Value = -379.5
$ws.Range("D38").Value = 948.65868
$ws.Range("C39").Value = 16117.76
$ws.Range("D39").Value = 444.55206

$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = 145
$ws.Range("C2").Value = -5.9091437
$ws.Range("B3").Value = 120
$ws.Range("C3").Value = -21299.687
$ws.Range("B4").Value = -41971.202
$ws.Range("C4").Value = -4.413347
$ws.Range("B5").Value = -2257.1266
$ws.Range("C5").Value = -9597.1353
$ws.Range("B6").Value = 2237.4494
$ws.Range("C6").Value = 88.675713
$ws.Range("B7").Value = 3370.9091
$ws.Range("C7").Value = -4995.4433
$ws.Range("C8").Value = 25774.432
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 14851.717
$ws.Range("B10").Value = -4203.75
$ws.Range("C10").Value = -6.247246
$ws.Range("B11").Value = 4403.75
$ws.Range("C11").Value = 192.84686
$ws.Range("B12").Value = 12265.587
$ws.Range("C12").Value = 87.792779
$ws.Range("B13").Value = 2223.0177
$ws.Range("C13").Value = 95.55015299999999
$ws.Range("B14").Value = 5699.6799
$ws.Range("C14").Value = 98.48417000000001
$ws.Range("B15").Value = 30176.344
$ws.Range("C15").Value = 72.930785
$ws.Range("B16").Value = -1630.9927
$ws.Range("C16").Value = 44.546653
$ws.Range("B17").Value = 55336.692
$ws.Range("C17").Value = 63.499222
$ws.Range("B18").Value = 686.44068
$ws.Range("C18").Value = 75.19922200000001
$ws.Range("C19").Value = 227.71689
$ws.Range("B20").Value = -2573.0698
$ws.Range("C20").Value = 73.849222
$ws.Range("B21").Value = -848.23583
$ws.Range("C21").Value = 86.469222
$ws.Range("B22").Value = -444.55206
$ws.Range("C22").Value = 27.324368
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = -1068.6833
$ws.Range("B24").Value = 100.42285
$ws.Range("C24").Value = 94.81822200000001
$ws.Range("B25").Value = -182004.46
$ws.Range("C25").Value = 25.826653
